$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Harvard case classification: recomputed "_old" baseline columns and
# reshuffled the average_doctor / average_doctor_old pair (a new
# average_doctor_old column was inserted before the prior average_doctor
# column, pushing its data into BQ and swapping the header labels).

$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"
$ws.Range("E4").Value = 0.385
$ws.Range("F4").Value = 0.082
$ws.Range("G4").Value = 0.287
$ws.Range("N4").Value = 0.404
$ws.Range("O4").Value = 0.062
$ws.Range("P4").Value = 0.25
$ws.Range("W4").Value = 0.217
$ws.Range("X4").Value = 0.101
$ws.Range("Y4").Value = 0.318
$ws.Range("AI4").Value = 0.197
$ws.Range("AJ4").Value = 0.064
$ws.Range("AK4").Value = 0.254
$ws.Range("AU4").Value = 0.148
$ws.Range("AV4").Value = 0.026
$ws.Range("AW4").Value = 0.161
$ws.Range("BA4").Value = 1.949
$ws.Range("BB4").Value = 0.169
$ws.Range("BC4").Value = 0.411
$ws.Range("BG4").Value = 0.73
$ws.Range("BH4").Value = 0.142
$ws.Range("BI4").Value = 0.377
$ws.Range("BM4").Value = 0.6820000000000001
$ws.Range("BN4").Value = 0.091
$ws.Range("BO4").Value = 0.302
$ws.Range("BP4").Value = 0.65
$ws.Range("BQ4").Value = 0.656
$ws.Range("E5").Value = 0.497
$ws.Range("F5").Value = 0.1
$ws.Range("G5").Value = 0.316
$ws.Range("N5").Value = 0.762
$ws.Range("O5").Value = 0.076
$ws.Range("P5").Value = 0.275
$ws.Range("W5").Value = 0.225
$ws.Range("X5").Value = 0.113
$ws.Range("Y5").Value = 0.337
$ws.Range("AI5").Value = 0.235
$ws.Range("AJ5").Value = 0.093
$ws.Range("AK5").Value = 0.305
$ws.Range("AU5").Value = 0.304
$ws.Range("AV5").Value = 0.098
$ws.Range("AW5").Value = 0.312
$ws.Range("BA5").Value = 1.367
$ws.Range("BB5").Value = 0.08599999999999999
$ws.Range("BC5").Value = 0.293
$ws.Range("BG5").Value = 0.41
$ws.Range("BI5").Value = 0.229
$ws.Range("BM5").Value = 0.574
$ws.Range("BP5").Value = 0.456
$ws.Range("BQ5").Value = 0.454
$ws.Range("E6").Value = 0.434
$ws.Range("N6").Value = 0.528
$ws.Range("W6").Value = 0.221
$ws.Range("AI6").Value = 0.214
$ws.Range("AU6").Value = 0.199
$ws.Range("BA6").Value = 1.595
$ws.Range("BG6").Value = 0.525
$ws.Range("BM6").Value = 0.623
$ws.Range("BP6").Value = 0.532
$ws.Range("BQ6").Value = 0.533
$ws.Range("E7").Value = 0.47
$ws.Range("N7").Value = 0.647
$ws.Range("W7").Value = 0.223
$ws.Range("AI7").Value = 0.226
$ws.Range("AU7").Value = 0.251
$ws.Range("BA7").Value = 1.448
$ws.Range("BG7").Value = 0.449
$ws.Range("BM7").Value = 0.593
$ws.Range("BP7").Value = 0.483
$ws.Range("BQ7").Value = 0.482
$ws.Range("E8").Value = 0.529
$ws.Range("F8").Value = 0.127
$ws.Range("G8").Value = 0.357
$ws.Range("N8").Value = 0.763
$ws.Range("O8").Value = 0.061
$ws.Range("P8").Value = 0.246
$ws.Range("W8").Value = 0.217
$ws.Range("X8").Value = 0.105
$ws.Range("Y8").Value = 0.325
$ws.Range("AI8").Value = 0.215
$ws.Range("AJ8").Value = 0.092
$ws.Range("AK8").Value = 0.303
$ws.Range("AU8").Value = 0.24
$ws.Range("AV8").Value = 0.075
$ws.Range("AW8").Value = 0.273
$ws.Range("BA8").Value = 1.716
$ws.Range("BB8").Value = 0.135
$ws.Range("BC8").Value = 0.367
$ws.Range("BG8").Value = 0.5649999999999999
$ws.Range("BH8").Value = 0.105
$ws.Range("BI8").Value = 0.324
$ws.Range("BM8").Value = 0.702
$ws.Range("BN8").Value = 0.07199999999999999
$ws.Range("BO8").Value = 0.268
$ws.Range("BP8").Value = 0.572
$ws.Range("BQ8").Value = 0.584
$ws.Range("E9").Value = 0.46
$ws.Range("F9").Value = 0.248
$ws.Range("G9").Value = 0.498
$ws.Range("N9").Value = 0.64
$ws.Range("O9").Value = 0.23
$ws.Range("P9").Value = 0.48
$ws.Range("W9").Value = 0.12
$ws.Range("X9").Value = 0.106
$ws.Range("Y9").Value = 0.325
$ws.Range("AI9").Value = 0.12
$ws.Range("AJ9").Value = 0.106
$ws.Range("AK9").Value = 0.325
$ws.Range("BA9").Value = 1.66
$ws.Range("BB9").Value = 0.244
$ws.Range("BC9").Value = 0.494
$ws.Range("BG9").Value = 0.6
$ws.Range("BH9").Value = 0.24
$ws.Range("BI9").Value = 0.49
$ws.Range("BM9").Value = 0.64
$ws.Range("BN9").Value = 0.23
$ws.Range("BO9").Value = 0.48
$ws.Range("BP9").Value = 0.553
$ws.Range("BQ9").Value = 0.554
$ws.Range("E10").Value = 0.58
$ws.Range("F10").Value = 0.244
$ws.Range("G10").Value = 0.494
$ws.Range("N10").Value = 0.84
$ws.Range("O10").Value = 0.134
$ws.Range("P10").Value = 0.367
$ws.Range("W10").Value = 0.26
$ws.Range("X10").Value = 0.192
$ws.Range("Y10").Value = 0.439
$ws.Range("AI10").Value = 0.24
$ws.Range("AJ10").Value = 0.182
$ws.Range("AK10").Value = 0.427
$ws.Range("AU10").Value = 0.22
$ws.Range("AV10").Value = 0.172
$ws.Range("AW10").Value = 0.414
$ws.Range("BA10").Value = 1.98
$ws.Range("BB10").Value = 0.25
$ws.Range("BC10").Value = 0.5
$ws.Range("BG10").Value = 0.64
$ws.Range("BH10").Value = 0.23
$ws.Range("BI10").Value = 0.48
$ws.Range("BM10").Value = 0.86
$ws.Range("BN10").Value = 0.12
$ws.Range("BO10").Value = 0.347
$ws.Range("BP10").Value = 0.66
$ws.Range("BQ10").Value = 0.6879999999999999
$ws.Range("E11").Value = 0.6
$ws.Range("F11").Value = 0.24
$ws.Range("G11").Value = 0.49
$ws.Range("N11").Value = 0.88
$ws.Range("O11").Value = 0.106
$ws.Range("P11").Value = 0.325
$ws.Range("W11").Value = 0.26
$ws.Range("X11").Value = 0.192
$ws.Range("Y11").Value = 0.439
$ws.Range("AI11").Value = 0.24
$ws.Range("AJ11").Value = 0.182
$ws.Range("AK11").Value = 0.427
$ws.Range("AU11").Value = 0.34
$ws.Range("AV11").Value = 0.224
$ws.Range("AW11").Value = 0.474
$ws.Range("BA11").Value = 1.98
$ws.Range("BB11").Value = 0.25
$ws.Range("BC11").Value = 0.5
$ws.Range("BG11").Value = 0.64
$ws.Range("BH11").Value = 0.23
$ws.Range("BI11").Value = 0.48
$ws.Range("BM11").Value = 0.86
$ws.Range("BN11").Value = 0.12
$ws.Range("BO11").Value = 0.347
$ws.Range("BP11").Value = 0.66
$ws.Range("BQ11").Value = 0.6879999999999999
$ws.Range("E12").Value = 1.4
$ws.Range("F12").Value = 0.64
$ws.Range("G12").Value = 0.8
$ws.Range("N12").Value = 1.652
$ws.Range("O12").Value = 1.618
$ws.Range("P12").Value = 1.272
$ws.Range("W12").Value = 1.846
$ws.Range("X12").Value = 0.746
$ws.Range("Y12").Value = 0.863
$ws.Range("AI12").Value = 1.917
$ws.Range("AJ12").Value = 0.91
$ws.Range("AK12").Value = 0.954
$ws.Range("AU12").Value = 2.647
$ws.Range("AV12").Value = 1.758
$ws.Range("AW12").Value = 1.326
$ws.Range("BA12").Value = 3.619
$ws.Range("BB12").Value = 0.332
$ws.Range("BC12").Value = 0.576
$ws.Range("BG12").Value = 1.062
$ws.Range("BH12").Value = 0.059
$ws.Range("BI12").Value = 0.242
$ws.Range("BM12").Value = 1.349
$ws.Range("BN12").Value = 0.413
$ws.Range("BO12").Value = 0.643
$ws.Range("BP12").Value = 1.206
$ws.Range("BQ12").Value = 1.27
$ws.Range("E13").Value = 1.728
$ws.Range("F13").Value = 0.92
$ws.Range("G13").Value = 0.959
$ws.Range("N13").Value = 2.338
$ws.Range("O13").Value = 1.166
$ws.Range("P13").Value = 1.08
$ws.Range("W13").Value = 1.09
$ws.Range("X13").Value = 0.186
$ws.Range("Y13").Value = 0.431
$ws.Range("AI13").Value = 1.383
$ws.Range("AJ13").Value = 0.401
$ws.Range("AK13").Value = 0.633
$ws.Range("AU13").Value = 2.482
$ws.Range("AV13").Value = 1.307
$ws.Range("AW13").Value = 1.143
$ws.Range("BA13").Value = 2.519
$ws.Range("BB13").Value = 0.313
$ws.Range("BC13").Value = 0.5590000000000001
$ws.Range("BG13").Value = 0.621
$ws.Range("BH13").Value = 0.08699999999999999
$ws.Range("BI13").Value = 0.294
$ws.Range("BM13").Value = 0.999
$ws.Range("BN13").Value = 0.363
$ws.Range("BO13").Value = 0.603
$ws.Range("BP13").Value = 0.84
$ws.Range("BQ13").Value = 0.786

